$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.005723000769734084
$ws.Range("J2").Value = 0.005723000769734084
$ws.Range("M2").Value = 0.478362
$ws.Range("N2").Value = 1.435086
$ws.Range("O2").Value = 0.6533545125880439
$ws.Range("P2").Value = 0.6533545125880439
$ws.Range("Q2").Value = 0.134972070656
$ws.Range("R2").Value = 1.214748635904
$ws.Range("S2").Value = 0.003739148378450613
$ws.Range("T2").Value = 0.003739148378450613
$ws.Range("I3").Value = 0.005723000769734084
$ws.Range("J3").Value = 0.005723000769734084
$ws.Range("O3").Value = 0.2124690813384451
$ws.Range("P3").Value = 0.2124690813384451
$ws.Range("S3").Value = 0.001215960716044615
$ws.Range("T3").Value = 0.001215960716044615
$ws.Range("I4").Value = 0.005723000769734084
$ws.Range("J4").Value = 0.005723000769734084
$ws.Range("M4").Value = 0.09823900000000001
$ws.Range("N4").Value = 0.294717
$ws.Range("O4").Value = 0.134176406073511
$ws.Range("P4").Value = 0.1341764060735109
$ws.Range("Q4").Value = 0.02771859229866667
$ws.Range("R4").Value = 0.249467330688
$ws.Range("S4").Value = 0.0007678916752388564
$ws.Range("T4").Value = 0.0007678916752388562
$ws.Range("G5").Value = 30.199365
$ws.Range("H5").Value = 90.598095
$ws.Range("I5").Value = 0.6125398923302606
$ws.Range("J5").Value = 0.6125398923302606
$ws.Range("M5").Value = 0.478362
$ws.Range("N5").Value = 1.435086
$ws.Range("O5").Value = 0.6533545125880439
$ws.Range("P5").Value = 0.6533545125880439
$ws.Range("Q5").Value = 14.44622864013
$ws.Range("R5").Value = 130.01605776117
$ws.Range("S5").Value = 0.4002057027941703
$ws.Range("T5").Value = 0.4002057027941703
$ws.Range("G6").Value = 30.199365
$ws.Range("H6").Value = 90.598095
$ws.Range("I6").Value = 0.6125398923302606
$ws.Range("J6").Value = 0.6125398923302606
$ws.Range("O6").Value = 0.2124690813384451
$ws.Range("P6").Value = 0.2124690813384451
$ws.Range("Q6").Value = 4.69787361813
$ws.Range("R6").Value = 42.28086256317
$ws.Range("S6").Value = 0.1301457882065606
$ws.Range("T6").Value = 0.1301457882065606
$ws.Range("G7").Value = 30.199365
$ws.Range("H7").Value = 90.598095
$ws.Range("I7").Value = 0.6125398923302606
$ws.Range("J7").Value = 0.6125398923302606
$ws.Range("M7").Value = 0.09823900000000001
$ws.Range("N7").Value = 0.294717
$ws.Range("O7").Value = 0.134176406073511
$ws.Range("P7").Value = 0.1341764060735109
$ws.Range("Q7").Value = 2.966755418235
$ws.Range("R7").Value = 26.700798764115
$ws.Range("S7").Value = 0.08218840132952974
$ws.Range("T7").Value = 0.08218840132952972
$ws.Range("G8").Value = 18.820355
$ws.Range("H8").Value = 56.461065
$ws.Range("I8").Value = 0.3817371069000054
$ws.Range("J8").Value = 0.3817371069000054
$ws.Range("M8").Value = 0.478362
$ws.Range("N8").Value = 1.435086
$ws.Range("O8").Value = 0.6533545125880439
$ws.Range("P8").Value = 0.6533545125880439
$ws.Range("Q8").Value = 9.002942658509999
$ws.Range("R8").Value = 81.02648392659
$ws.Range("S8").Value = 0.249409661415423
$ws.Range("T8").Value = 0.249409661415423
$ws.Range("G9").Value = 18.820355
$ws.Range("H9").Value = 56.461065
$ws.Range("I9").Value = 0.3817371069000054
$ws.Range("J9").Value = 0.3817371069000054
$ws.Range("O9").Value = 0.2124690813384451
$ws.Range("P9").Value = 0.2124690813384451
$ws.Range("Q9").Value = 2.92773206451
$ws.Range("R9").Value = 26.34958858059
$ws.Range("S9").Value = 0.08110733241583998
$ws.Range("T9").Value = 0.08110733241583996
$ws.Range("G10").Value = 18.820355
$ws.Range("H10").Value = 56.461065
$ws.Range("I10").Value = 0.3817371069000054
$ws.Range("J10").Value = 0.3817371069000054
$ws.Range("M10").Value = 0.09823900000000001
$ws.Range("N10").Value = 0.294717
$ws.Range("O10").Value = 0.134176406073511
$ws.Range("P10").Value = 0.1341764060735109
$ws.Range("Q10").Value = 1.848892854845
$ws.Range("R10").Value = 16.640035693605
$ws.Range("S10").Value = 0.05122011306874239
$ws.Range("T10").Value = 0.05122011306874238
